$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("PE","LIMA METROPOLITANA",107563,189192,296755,11989,1466082,44079),
    @("PE","CALLAO",7113,21293,28406,1600,109819,44079),
    @("PE","LA LIBERTAD",3371,20150,23521,2091,119408,44079),
    @("PE","AREQUIPA",4500,30285,34785,1199,200867,44079),
    @("PE","MOQUEGUA",430,9231,9661,237,59357,44079),
    @("PE","CUSCO",5473,10962,16435,342,90499,44079),
    @("PE","LAMBAYEQUE",2430,20433,22863,1606,94828,44079),
    @("PE","PIURA",1024,29398,30422,1848,117611,44079),
    @("PE","AMAZONAS",993,10485,11478,192,56441,44079),
    @("PE","ICA",3143,19200,22343,1484,108999,44079),
    @("PE","JUNIN",2114,13773,15887,678,90960,44079),
    @("PE","ANCASH",3574,15659,19233,1189,99239,44079),
    @("PE","CAJAMARCA",2675,13650,16325,400,88048,44079),
    @("PE","HUANUCO",1303,11822,13125,344,65688,44079),
    @("PE","APURIMAC",394,2854,3248,69,39565,44079),
    @("PE","PUNO",650,8918,9568,250,60392,44079),
    @("PE","TACNA",1451,7856,9307,159,55907,44079),
    @("PE","LORETO",3160,13149,16309,927,51673,44079),
    @("PE","SAN MARTIN",2444,11940,14384,665,65605,44079),
    @("PE","AYACUCHO",2227,6705,8932,273,45715,44079),
    @("PE","HUANCAVELICA",719,4447,5166,88,35731,44079),
    @("PE","PASCO",263,3839,4102,91,32889,44079),
    @("PE","TUMBES",680,5268,5948,283,26780,44079),
    @("PE","UCAYALI",1004,13315,14319,302,46065,44079),
    @("PE","MADRE DE DIOS",900,5518,6418,131,30750,44079),
    @("PE","LIMA REGIÓN",5513,19249,24762,1250,94517,44079),
    @("PE","LIMA METROPOLITANA",108292,190622,298914,12054,1475932,44080),
    @("PE","CALLAO",7140,21467,28607,1608,110458,44080),
    @("PE","LA LIBERTAD",3389,20345,23734,2096,120049,44080),
    @("PE","AREQUIPA",4504,30649,35153,1205,202341,44080),
    @("PE","MOQUEGUA",430,9339,9769,238,59632,44080),
    @("PE","CUSCO",5509,11171,16680,349,91400,44080),
    @("PE","LAMBAYEQUE",2450,20578,23028,1610,95361,44080),
    @("PE","PIURA",1025,29647,30672,1853,118451,44080),
    @("PE","AMAZONAS",993,10618,11611,194,57494,44080),
    @("PE","ICA",3168,19497,22665,1486,110224,44080),
    @("PE","JUNIN",2114,13963,16077,682,91820,44080),
    @("PE","ANCASH",3603,15791,19394,1197,99905,44080),
    @("PE","CAJAMARCA",2692,13783,16475,403,88728,44080),
    @("PE","HUANUCO",1303,11982,13285,346,66202,44080),
    @("PE","APURIMAC",394,2960,3354,69,41182,44080),
    @("PE","PUNO",650,9156,9806,256,61115,44080),
    @("PE","TACNA",1458,8061,9519,163,56412,44080),
    @("PE","LORETO",3169,13315,16484,929,52075,44080),
    @("PE","SAN MARTIN",2469,12120,14589,669,66085,44080),
    @("PE","AYACUCHO",2235,6837,9072,276,46331,44080),
    @("PE","HUANCAVELICA",719,4479,5198,89,35916,44080),
    @("PE","PASCO",270,3882,4152,92,33240,44080),
    @("PE","TUMBES",680,5300,5980,283,26916,44080),
    @("PE","UCAYALI",1004,13385,14389,302,46280,44080),
    @("PE","MADRE DE DIOS",900,5592,6492,131,30966,44080),
    @("PE","LIMA REGIÓN",5525,19353,24878,1258,95065,44080)
)

$startRow = 4095
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

$endRow = $startRow + $rows.Count - 1
$ws.Cells.Item(4094, 8).Copy()
$ws.Range("H" + $startRow + ":H" + $endRow).PasteSpecial(-4122)
